# G1_LP_A.xlsx — "added new scale 1-5 rating of confidence for test phase"
#
# Updates a batch of shared-string cell values on Tabelle1 (two-letter
# syllable codes drawn from the workbook's existing shared-strings table),
# moves the active-cell selection, and turns off multithreaded/concurrent
# calculation to match the commit's workbook-level calcPr tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Disable concurrent (multi-threaded) calculation — mirrors
# <calcPr concurrentCalc="0"/> in workbook.xml.
$excel.Application.MultiThreadedCalculation.Enabled = $false

# Row 2
$ws.Range("F2").Value = "fi"

# Row 4
$ws.Range("B4").Value = "lo"
$ws.Range("C4").Value = "ba"
$ws.Range("E4").Value = "lo"
$ws.Range("F4").Value = "pe"
$ws.Range("G4").Value = "fi"

# Row 5
$ws.Range("C5").Value = "fu"
$ws.Range("D5").Value = "lo"
$ws.Range("E5").Value = "to"
$ws.Range("F5").Value = "ko"
$ws.Range("G5").Value = "ba"

# Row 8
$ws.Range("F8").Value = "ba"

# Row 9
$ws.Range("B9").Value = "se"
$ws.Range("C9").Value = "lo"

# Row 12
$ws.Range("D12").Value = "ba"
$ws.Range("E12").Value = "fu"

# Row 16
$ws.Range("B16").Value = "to"
$ws.Range("C16").Value = "se"
$ws.Range("D16").Value = "pe"
$ws.Range("E16").Value = "fu"
$ws.Range("F16").Value = "to"
$ws.Range("G16").Value = "pe"

# Row 18
$ws.Range("D18").Value = "fu"
$ws.Range("E18").Value = "to"
$ws.Range("F18").Value = "ba"
$ws.Range("G18").Value = "fu"

# Row 24
$ws.Range("F24").Value = "pe"
$ws.Range("G24").Value = "fu"

# Row 26
$ws.Range("D26").Value = "fu"

# Row 28
$ws.Range("E28").Value = "di"
$ws.Range("F28").Value = "ni"

# Row 30
$ws.Range("F30").Value = "fu"

# Move the active selection on Tabelle1 from N4 to D33.
$ws.Range("D33").Select() | Out-Null
